$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D6: replace the old "https://exemple.com" text with "https://youtube.com"
# and apply the same "Lien hypertexte" (hyperlink) style used by D7/D8.
$ws.Range("D6").Value = "https://youtube.com"
$ws.Range("D6").Style = "Lien hypertexte"

# Update the active selection to D7 (as reflected in the saved workbook view)
$ws.Range("D7").Select()
